$wb = $excel.ActiveWorkbook

# --- 1. Update "Last Updated" timestamp on the Metadata sheet -------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value2 = "05 Nov 2025, 01:32 PM"

# --- 2. Insert a new leading stock ("CAPTRU-RE1") at the top of the       --
#        "Stock List" sheet, pushing every existing row down by one and   --
#        dropping the last row off the bottom of the list. -----------------
$ws = $wb.Worksheets.Item("Stock List")

$lastRow = 76
$firstDataRow = 2

# Capture the existing B/C/D/E/H values for every data row before
# overwriting anything, so the shift-down can be performed safely.
$oldB = @{}
$oldC = @{}
$oldD = @{}
$oldE = @{}
$oldH = @{}

for ($r = $firstDataRow; $r -le ($lastRow - 1); $r++) {
    $oldB[$r] = $ws.Cells.Item($r, 2).Value2
    $oldC[$r] = $ws.Cells.Item($r, 3).Value2
    $oldD[$r] = $ws.Cells.Item($r, 4).Value2
    $oldE[$r] = $ws.Cells.Item($r, 5).Value2
    $oldH[$r] = $ws.Cells.Item($r, 8).Value2
}

# Shift rows 2..75 down into rows 3..76 (row 76's prior contents are
# discarded, as the list keeps a fixed number of rows).
for ($r = ($lastRow - 1); $r -ge $firstDataRow; $r--) {
    $dest = $r + 1
    $ws.Cells.Item($dest, 2).Value2 = $oldB[$r]
    $ws.Cells.Item($dest, 3).Value2 = $oldC[$r]
    $ws.Cells.Item($dest, 4).Value2 = $oldD[$r]
    $ws.Cells.Item($dest, 5).Value2 = $oldE[$r]
    $ws.Cells.Item($dest, 8).Value2 = $oldH[$r]
}

# Write the brand-new top entry into row 2.
$ws.Cells.Item(2, 2).Value2 = "CAPTRU-RE1"
$ws.Cells.Item(2, 3).Value2 = "CAPTRU-RE1"
$ws.Cells.Item(2, 4).Value2 = 5.67
$ws.Cells.Item(2, 5).Value2 = -11.9565
$ws.Cells.Item(2, 8).Value2 = 0
